$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "['MEC-3B-Usin. CNC', -, -, -]"
$ws.Range("D2").Value = "-"

$ws.Range("B3").Value = "[-, 'MEC-2B-Ajustagem', -, -]"
$ws.Range("D3").Value = "[-, 'MEC-3B-Usin. CNC', -, -]"

$ws.Range("B4").Value = "[-, 'MEC-2B-Ajustagem', -, -]"
$ws.Range("D4").Value = "[-, 'MEC-3B-Usin. CNC', -, -]"
$ws.Range("F4").Value = "[-, -, -, 'MEC-2B-Ajustagem']"

$ws.Range("B6").Value = "[-, 'MEC-2B-Ajustagem', -, -]"
$ws.Range("D6").Value = "-"

$ws.Range("B7").Value = "-"

$ws.Range("C8").Value = "[-, -, -, 'MEC-3B-Usin. CNC']"

$ws.Range("B18").Value = "[-, -, -, 'MEC-2NB-Usin. CNC']"
$ws.Range("C18").Value = "[-, -, -, 'MEC-1NB-Ajustagem']"
$ws.Range("D18").Value = "-"
$ws.Range("E18").Value = "-"
$ws.Range("F18").Value = "[-, 'MEC-1NB-Ajustagem', 'MEC-1NA-Ajustagem', -]"

$ws.Range("B19").Value = "[-, -, -, 'MEC-2NB-Usin. CNC']"
$ws.Range("F19").Value = "[-, -, 'MEC-1NA-Ajustagem', -]"

$ws.Range("B20").Value = "-"
$ws.Range("C20").Value = "[-, -, -, 'MEC-1NB-Ajustagem']"
$ws.Range("F20").Value = "[-, -, 'MEC-1NA-Ajustagem', -]"

$ws.Range("B21").Value = "['MEC-1NB-Ajustagem', -, -, 'MEC-2NB-Usin. CNC']"
$ws.Range("D21").Value = "[-, -, 'MEC-2NB-Usin. CNC', -]"
$ws.Range("F21").Value = "[-, -, 'MEC-1NA-Ajustagem', -]"
